$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 3.5
$ws.Range("O2").Value = 1.36
$ws.Range("S2").Value = 3.8
$ws.Range("T2").Value = 1.86
$ws.Range("U2").Value = 2.06
$ws.Range("G6").Value = 9
$ws.Range("Y6").Value = 1000
$ws.Range("Z6").Value = 1000
$ws.Range("AA6").Value = 1000
$ws.Range("AC6").Value = 1000
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 1000
$ws.Range("AO6").Value = 1000
$ws.Range("O7").Value = 1.16
$ws.Range("P7").Value = 2.62
$ws.Range("Q7").Value = 1.49
$ws.Range("R7").Value = 1.66
$ws.Range("S7").Value = 2.2
$ws.Range("V7").Value = 1.79
$ws.Range("X7").Value = 36
$ws.Range("Y7").Value = 19.5
$ws.Range("Z7").Value = 22
$ws.Range("AA7").Value = 34
$ws.Range("AB7").Value = 25
$ws.Range("AC7").Value = 13
$ws.Range("AF7").Value = 36
$ws.Range("AG7").Value = 18.5
$ws.Range("AK7").Value = 38
$ws.Range("AL7").Value = 40
$ws.Range("AM7").Value = 65
$ws.Range("AN7").Value = 24
$ws.Range("I8").Value = 3.9
$ws.Range("J8").Value = 3.75
$ws.Range("T8").Value = 1.53
$ws.Range("V8").Value = 1.38
$ws.Range("W8").Value = 1.72
$ws.Range("Z8").Value = 32
$ws.Range("AC8").Value = 12
$ws.Range("AD8").Value = 17.5
$ws.Range("AI8").Value = 44
$ws.Range("AL8").Value = 34
$ws.Range("I9").Value = 1.75
$ws.Range("Q9").Value = 1.77
$ws.Range("S9").Value = 2.88
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 2.32
$ws.Range("X9").Value = 1000
$ws.Range("Y9").Value = 1000
$ws.Range("Z9").Value = 1000
$ws.Range("AA9").Value = 1000
$ws.Range("AD9").Value = 1000
$ws.Range("AE9").Value = 1000
$ws.Range("AO9").Value = 1000
$ws.Range("T10").Value = 1.6
$ws.Range("H11").Value = 1.63
$ws.Range("Q11").Value = 1.61
$ws.Range("G12").Value = 1.36
$ws.Range("H12").Value = 9.6
$ws.Range("K12").Value = 7.4
$ws.Range("N12").Value = 5.7
$ws.Range("T12").Value = 1.93
$ws.Range("AB12").Value = 13.5
$ws.Range("AC12").Value = 18.5
$ws.Range("AJ12").Value = 12.5
$ws.Range("AN12").Value = 4.9
